$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("J2").Value = 688307.5510076137
$ws.Range("K2").Value = 633186.9419317616
$ws.Range("L2").Value = 775035.2065316344
$ws.Range("M2").Value = 712420.0651001015
$ws.Range("J3").Value = 1260058.578931904
$ws.Range("K3").Value = 1282281.012875688
$ws.Range("L3").Value = 1413255.438083681
$ws.Range("M3").Value = 1422239.432561247
$ws.Range("J4").Value = 738589.5847541089
$ws.Range("K4").Value = 751950.1743426092
$ws.Range("L4").Value = 829544.8904273655
$ws.Range("M4").Value = 835045.6171247229
$ws.Range("J5").Value = 1331409.121111143
$ws.Range("K5").Value = 1239039.106639018
$ws.Range("L5").Value = 1501156.812972522
$ws.Range("M5").Value = 1394118.126076698
$ws.Range("J6").Value = 2243679.95666909
$ws.Range("K6").Value = 2157428.560842621
$ws.Range("L6").Value = 2534705.711240039
$ws.Range("M6").Value = 2423305.452334199
$ws.Range("J7").Value = 2954652.168463362
$ws.Range("K7").Value = 2783128.465630074
$ws.Range("L7").Value = 3330808.868291364
$ws.Range("M7").Value = 3126779.749041206
$ws.Range("J8").Value = 20540821.20948786
$ws.Range("K8").Value = 20751118.86907979
$ws.Range("L8").Value = 23095889.76559405
$ws.Range("M8").Value = 23085392.29531628
$ws.Range("J9").Value = 16441623.00424981
$ws.Range("K9").Value = 16520789.41145229
$ws.Range("L9").Value = 18520491.03351302
$ws.Range("M9").Value = 18420012.97074377
$ws.Range("J10").Value = 1221222.462646787
$ws.Range("K10").Value = 1189871.359363376
$ws.Range("L10").Value = 1379593.369667514
$ws.Range("M10").Value = 1334556.713838103
$ws.Range("J11").Value = 7646822.111145169
$ws.Range("K11").Value = 7589009.35114764
$ws.Range("L11").Value = 8624957.234830759
$ws.Range("M11").Value = 8482619.356652685
$ws.Range("J12").Value = 3330162.920451522
$ws.Range("K12").Value = 3162761.516548321
$ws.Range("L12").Value = 3758564.062899727
$ws.Range("M12").Value = 3554142.561959549
$ws.Range("J13").Value = 2260558.937626127
$ws.Range("K13").Value = 2361576.390307112
$ws.Range("L13").Value = 2560065.788831091
$ws.Range("M13").Value = 2635201.488182747
$ws.Range("J14").Value = 6553021.814271568
$ws.Range("K14").Value = 6594814.295873213
$ws.Range("L14").Value = 7419415.757226863
$ws.Range("M14").Value = 7386339.185825613
